$wb = $excel.ActiveWorkbook

# Sheet "展览" (1st sheet) - column F "想去人数" updates
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 14982
$ws1.Range("F3").Value = 18931
$ws1.Range("F5").Value = 136
$ws1.Range("F14").Value = 148
$ws1.Range("F15").Value = 216
$ws1.Range("F17").Value = 1448
$ws1.Range("F20").Value = 95
$ws1.Range("F22").Value = 7877
$ws1.Range("F27").Value = 1238
$ws1.Range("F29").Value = 6029
$ws1.Range("F30").Value = 113
$ws1.Range("F34").Value = 276
$ws1.Range("F35").Value = 5396
$ws1.Range("F36").Value = 18

# Sheet "全部类型" (4th sheet) - column F "想去人数" updates
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 14982
$ws4.Range("F3").Value = 18931
$ws4.Range("F5").Value = 136
$ws4.Range("F14").Value = 148
$ws4.Range("F15").Value = 216
$ws4.Range("F17").Value = 1448
$ws4.Range("F21").Value = 95
$ws4.Range("F23").Value = 7877
$ws4.Range("F28").Value = 1238
$ws4.Range("F32").Value = 6029
$ws4.Range("F33").Value = 113
$ws4.Range("F37").Value = 276
$ws4.Range("F38").Value = 5396
$ws4.Range("F39").Value = 18
